# Generate Report for Handoff
# Adds a new row (row 3) to the Overview, zh-cn, and de-de sheets for the
# newly handed-off file "d8bc7293-8054-49a3-940d-ff0a336231c1oooo...md",
# mirroring the existing row 2 pattern, and widens a couple of columns that
# now need to fit the longer generated filenames.

$wb = $excel.ActiveWorkbook

$newUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bd64f31bb88ba742c99a63f1022599d417151ec/e2e/d8bc7293-8054-49a3-940d-ff0a336231c1ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "d8bc7293-8054-49a3-940d-ff0a336231c1ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("D3").Interior.ColorIndex = -4142
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"

$wsOverview.Range("G3").Value = "2016-08-20 18:40:19"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", "e2e\d8bc7293-8054-49a3-940d-ff0a336231c1ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"

$wsZhCn.Range("G3").Value = "d8bc7293-8054-49a3-940d-ff0a336231c1oooooooooooooooooooooooooooooooooooooooo.5b6e63e79a8165510c67d3e0e73a2dd6a04ad5e6.zh-cn.xlf"

$wsZhCn.Range("H3").Value = "2016-08-20 18:40:15"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("I3").Interior.ColorIndex = -4142
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("J3").Interior.ColorIndex = -4142

$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("L3").Interior.ColorIndex = -4142
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("N3").Interior.ColorIndex = -4142
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Range("P3").Interior.ColorIndex = -4142

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newUrl, "", "", "d8bc7293-8054-49a3-940d-ff0a336231c1ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"

$wsDeDe.Range("G3").Value = "d8bc7293-8054-49a3-940d-ff0a336231c1oooooooooooooooooooooooooooooooooooooooo.5b6e63e79a8165510c67d3e0e73a2dd6a04ad5e6.de-de.xlf"

$wsDeDe.Range("H3").Value = "2016-08-20 18:40:19"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("I3").Interior.ColorIndex = -4142
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("J3").Interior.ColorIndex = -4142

$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("L3").Interior.ColorIndex = -4142
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("N3").Interior.ColorIndex = -4142
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Range("P3").Interior.ColorIndex = -4142

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newUrl, "", "", "d8bc7293-8054-49a3-940d-ff0a336231c1ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md")

$wsDeDe.Columns.Item(3).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# Extend the tables (ListObjects) so the new row is included in each
# table's range / autofilter, matching the new A1:*3 extents.
# ---------------------------------------------------------------------
foreach ($wsName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($wsName)
    if ($ws.ListObjects.Count -gt 0) {
        $lo = $ws.ListObjects.Item(1)
        $lo.Resize($lo.Range.Resize(3))
    }
}
